# This workbook's data rows (2-23) are being reshuffled: the set of values in
# columns D (Fecha), I (Calidad), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), N (Unidad de
# comercializacion), P (Precio $/Kg) and Q (Kg o Unidades) move between rows
# according to a fixed permutation (columns A,B,C,E,F,G,H,O,R stay the same
# for every row already, so they are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: key = destination row, value = source row (i.e. destination row
# receives the values that used to live in the source row).
$map = @{
    2  = 23
    3  = 18
    4  = 5
    5  = 17
    6  = 13
    7  = 11
    8  = 19
    9  = 10
    10 = 4
    11 = 16
    12 = 9
    13 = 14
    14 = 2
    15 = 12
    16 = 22
    17 = 6
    18 = 3
    19 = 8
    20 = 15
    21 = 7
    22 = 20
    23 = 21
}

# Snapshot the current ("before") values of every row first, since rows will
# be overwritten in place and some rows are part of multi-row permutation
# cycles.
$snapshot = @{}
for ($r = 2; $r -le 23; $r++) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        I = $ws.Cells.Item($r, 9).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
    }
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $src = $snapshot[$srcRow]

    $ws.Cells.Item($destRow, 4).Value2  = $src.D
    $ws.Cells.Item($destRow, 9).Value   = $src.I
    $ws.Cells.Item($destRow, 10).Value2 = $src.J
    $ws.Cells.Item($destRow, 11).Value2 = $src.K
    $ws.Cells.Item($destRow, 12).Value2 = $src.L
    $ws.Cells.Item($destRow, 13).Value2 = $src.M
    $ws.Cells.Item($destRow, 14).Value  = $src.N
    $ws.Cells.Item($destRow, 16).Value2 = $src.P
    $ws.Cells.Item($destRow, 17).Value2 = $src.Q
}
